$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-write the "Name und Vorname..." placeholder text with the corrected
# wording (space after "Nr." before "#idNr", single blank-line/newline
# before "Betrieblicher ..."). Re-assigning the value causes the shared
# string to be re-created (moving it to the end of the shared-strings
# table once the old, now-unreferenced entry is dropped on save).
$ws.Range("A1").Value = "Name und Vorname des Auszubildenen #idName`n#idYear. Ausbildungsjahr`nAusbildungsnachweis Nr. #idNr`nFür die Woche vom #idFirstDate bis #idLastDate.`nBetrieblicher Funktionsberreich: #idDepartment"

# Update the sheet's selection to span the header row.
$ws.Range("A1:F1").Select()
